# Conserto do erro com o rótulo da coluna 2050 nas tabelas e
# retirada das linhas com total das tabelas

$wb = $excel.ActiveWorkbook

# Sheets that have the "2050" single-year header in E1, a "Total" row
# in row 13 that must be removed (rows 1-12 keep the data).
$yearSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $yearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "'2050"
    $ws.Rows.Item(13).Delete()
}

# Sheet with decade ranges header ("2015-2030"/"2031-2040") whose E1
# must become "2041-2050"; also has a "Total" row in row 13 to remove.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("E1").Value = "'2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet with only the mislabeled E1 header to fix, no Total row present.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws5.Range("E1").Value = "'2050"

# Sheet with a "Total" row (row 4) that must be removed.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
